# chore: update Sheets via scheduled runner
# Refreshes market-price derived columns (currentAveragePrice /
# currentAveragePriceNQ / currentAveragePriceHQ / LevePriceNQ / LevePriceHQ /
# LeveProfitNQ / LeveProfitHQ) for specific Leve rows across the ALC, ARM,
# BSM, CRP, CUL, LTW and WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 38
$ws.Range("H38").Value = 664.9231
$ws.Range("I38").Value = 201.75
$ws.Range("J38").Value = 1406
$ws.Range("K38").Value = 605.25
$ws.Range("L38").Value = 4218
$ws.Range("M38").Value = -233.25
$ws.Range("N38").Value = -4962

# Row 42
$ws.Range("H42").Value = 404
$ws.Range("I42").Value = 8
$ws.Range("J42").Value = 800
$ws.Range("K42").Value = 24
$ws.Range("L42").Value = 2400
$ws.Range("M42").Value = 206
$ws.Range("N42").Value = -2860

# Row 121
$ws.Range("H121").Value = 982.9167
$ws.Range("J121").Value = 982.9167
$ws.Range("L121").Value = 2948.7501
$ws.Range("N121").Value = -6442.7501

# Row 131
$ws.Range("H131").Value = 2887.8333
$ws.Range("I131").Value = 680.5
$ws.Range("K131").Value = 2041.5
$ws.Range("M131").Value = 2998.5

# Row 137
$ws.Range("H137").Value = 2493.7385
$ws.Range("I137").Value = 2367.1304
$ws.Range("J137").Value = 2800.2632
$ws.Range("K137").Value = 7101.3912
$ws.Range("L137").Value = 8400.7896
$ws.Range("M137").Value = -4551.3912
$ws.Range("N137").Value = -13500.7896

# Row 138
$ws.Range("H138").Value = 3614.795
$ws.Range("I138").Value = 1246.1305
$ws.Range("J138").Value = 7019.75
$ws.Range("K138").Value = 3738.3915
$ws.Range("L138").Value = 21059.25
$ws.Range("M138").Value = 1401.6085
$ws.Range("N138").Value = -31339.25

# Row 141
$ws.Range("H141").Value = 3741.0667
$ws.Range("I141").Value = 1110.091
$ws.Range("J141").Value = 10976.25
$ws.Range("K141").Value = 3330.273
$ws.Range("L141").Value = 32928.75
$ws.Range("M141").Value = 1849.727
$ws.Range("N141").Value = -43288.75

$ws = $wb.Worksheets.Item("ARM")
# Row 31
$ws.Range("H31").Value = 2900
$ws.Range("I31").Value = 2900
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 2900
$ws.Range("L31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -2606

# Row 32
$ws.Range("H32").Value = 4020054.8
$ws.Range("I32").Value = 3477.4429
$ws.Range("K32").Value = 3477.4429
$ws.Range("M32").Value = -3190.4429

# Row 109
$ws.Range("H109").Value = 48629.332
$ws.Range("J109").Value = 48629.332
$ws.Range("L109").Value = 48629.332
$ws.Range("N109").Value = -51403.332

$ws = $wb.Worksheets.Item("BSM")
# Row 102
$ws.Range("H102").Value = 8000
$ws.Range("I102").Value = 8000
$ws.Range("K102").Value = 8000
$ws.Range("M102").Value = -4755

# Row 134
$ws.Range("H134").Value = 2061853.8
$ws.Range("I134").Value = 2342.818
$ws.Range("J134").Value = 11123702
$ws.Range("K134").Value = 7028.454000000001
$ws.Range("L134").Value = 33371106
$ws.Range("M134").Value = -4493.454000000001
$ws.Range("N134").Value = -33376176

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1175.6833
$ws.Range("I31").Value = 884.2105
$ws.Range("J31").Value = 1310.7561
$ws.Range("K31").Value = 884.2105
$ws.Range("L31").Value = 1310.7561
$ws.Range("M31").Value = -589.2105
$ws.Range("N31").Value = -1900.7561

# Row 34
$ws.Range("H34").Value = 1175.6833
$ws.Range("I34").Value = 884.2105
$ws.Range("J34").Value = 1310.7561
$ws.Range("K34").Value = 884.2105
$ws.Range("L34").Value = 1310.7561
$ws.Range("M34").Value = -682.2105
$ws.Range("N34").Value = -1714.7561

# Row 58
$ws.Range("H58").Value = 28572300
$ws.Range("I58").Value = 41667440
$ws.Range("J58").Value = 1086.2727
$ws.Range("K58").Value = 41667440
$ws.Range("L58").Value = 1086.2727
$ws.Range("M58").Value = -41667237
$ws.Range("N58").Value = -1492.2727

# Row 136
$ws.Range("H136").Value = 28572300
$ws.Range("I136").Value = 41667440
$ws.Range("J136").Value = 1086.2727
$ws.Range("K136").Value = 125002320
$ws.Range("L136").Value = 3258.8181
$ws.Range("M136").Value = -124999770
$ws.Range("N136").Value = -8358.8181

$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 6059.524
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 6059.524
$ws.Range("K68").Value = 0
$ws.Range("L68").ClearContents()
$ws.Range("M68").Value = 18178.572
$ws.Range("N68").Value = -19800.572

# Row 71
$ws.Range("H71").Value = 6059.524
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 6059.524
$ws.Range("K71").Value = 0
$ws.Range("L71").ClearContents()
$ws.Range("M71").Value = 54535.716
$ws.Range("N71").Value = -62647.716

# Row 81
$ws.Range("H81").Value = 2899.889
$ws.Range("I81").Value = 699.6667
$ws.Range("J81").Value = 4000
$ws.Range("K81").Value = 2099.0001
$ws.Range("L81").Value = 12000
$ws.Range("M81").Value = -976.0001000000002
$ws.Range("N81").Value = -14246

# Row 84
$ws.Range("H84").Value = 2899.889
$ws.Range("I84").Value = 699.6667
$ws.Range("J84").Value = 4000
$ws.Range("K84").Value = 6297.0003
$ws.Range("L84").Value = 36000
$ws.Range("M84").Value = -681.0002999999997
$ws.Range("N84").Value = -47232

# Row 131
$ws.Range("H131").Value = 916.6083
$ws.Range("I131").Value = 391.5263
$ws.Range("J131").Value = 1044.5128
$ws.Range("K131").Value = 1174.5789
$ws.Range("L131").Value = 3133.5384
$ws.Range("M131").Value = 3865.4211
$ws.Range("N131").Value = -13213.5384

$ws = $wb.Worksheets.Item("LTW")
# Row 82
$ws.Range("H82").Value = 1563.909
$ws.Range("I82").Value = 1150.3334
$ws.Range("J82").Value = 2060.2
$ws.Range("K82").Value = 1150.3334
$ws.Range("L82").Value = 2060.2
$ws.Range("M82").Value = -789.3334
$ws.Range("N82").Value = -2782.2

# Row 85
$ws.Range("H85").Value = 1563.909
$ws.Range("I85").Value = 1150.3334
$ws.Range("J85").Value = 2060.2
$ws.Range("K85").Value = 1150.3334
$ws.Range("L85").Value = 2060.2
$ws.Range("M85").Value = 97.66660000000002
$ws.Range("N85").Value = -4556.2

$ws = $wb.Worksheets.Item("WVR")
# Row 8
$ws.Range("H8").Value = 2633.3333
$ws.Range("J8").Value = 3900
$ws.Range("L8").Value = 3900
$ws.Range("N8").Value = -4180

# Row 70
$ws.Range("H70").Value = 23188.637
$ws.Range("I70").Value = 19537.5
$ws.Range("J70").Value = 24000
$ws.Range("K70").Value = 19537.5
$ws.Range("L70").Value = 24000
$ws.Range("M70").Value = -19222.5
$ws.Range("N70").Value = -24630

# Row 73
$ws.Range("H73").Value = 23188.637
$ws.Range("I73").Value = 19537.5
$ws.Range("J73").Value = 24000
$ws.Range("K73").Value = 19537.5
$ws.Range("L73").Value = 24000
$ws.Range("M73").Value = -18445.5
$ws.Range("N73").Value = -26184

# Row 132
$ws.Range("H132").Value = 20746.373
$ws.Range("I132").Value = 29105.217
$ws.Range("J132").Value = 6688.3184
$ws.Range("K132").Value = 87315.651
$ws.Range("L132").Value = 20064.9552
$ws.Range("M132").Value = -84785.651
$ws.Range("N132").Value = -25124.9552
